# Generate Report for Archive
#
# 1) The localization status for the one tracked file flips from
#    "Ready for handoff" to "In Translation" -- update every cell that
#    shows it: Overview!E2:F2 (per-language status columns) and the
#    "Status" column (C2) on each language detail sheet.
# 2) The (now shorter) status text means the "Status" column no longer
#    needs to be as wide, so narrow it on all three sheets.
#    ColumnWidth is pixel-quantised by the host, so 12.5 is the
#    character-width input that lands on the closest achievable pixel
#    grid point to the target ~13.41 width.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newColumnWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
